# Weekly refresh of the "Hortaliza, Agrícola del Norte S.A. de Arica - Pepino dulce"
# subset: the daily price rows are re-pulled/re-aligned against the weekly series,
# which reshuffles the Fecha/Variedad/Calidad/Volumen/Precio/Unidad/Origen/Clasificación
# figures across several rows. Row 1 (headers) and rows 3, 4, 11, 12 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> new D,H,I,J,K,L,M,N,O,P,Q values
$rows = @{
    2  = @(44412, 'Cultivar IV Región', 'Primera', 150, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18)
    5  = @(44221, 'Cultivar XV región', 'Primera', 140, 5000, 6000, 5500, '$/caja 10 kilos', 'Región de Arica y Parinacota', 550, 10)
    6  = @(44211, 'Cultivar XV región', 'Segunda', 140, 4500, 5000, 4750, '$/caja 10 kilos', 'Región de Arica y Parinacota', 475, 10)
    7  = @(44377, 'Cultivar IV Región', 'Primera', 100, 17000, 18000, 17600, '$/bandeja 18 kilos', 'Provincia de Limarí', 978, 18)
    8  = @(44435, 'Cultivar IV Región', 'Segunda', 100, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18)
    9  = @(44435, 'Cultivar IV Región', 'Tercera', 120, 14000, 15000, 14500, '$/bandeja 18 kilos', 'Provincia de Limarí', 806, 18)
    10 = @(44363, 'Cultivar IV Región', 'Primera', 140, 14000, 15000, 14500, '$/bandeja 18 kilos', 'Provincia de Limarí', 806, 18)
    13 = @(44526, 'Cultivar XV región', 'Primera', 100, 5000, 5500, 5250, '$/caja 10 kilos', 'Región de Arica y Parinacota', 525, 10)
    14 = @(44526, 'Cultivar XV región', 'Segunda', 100, 4000, 4500, 4250, '$/caja 10 kilos', 'Región de Arica y Parinacota', 425, 10)
    15 = @(44526, 'Cultivar XV región', 'Tercera', 120, 3000, 3500, 3250, '$/caja 10 kilos', 'Región de Arica y Parinacota', 325, 10)
    16 = @(44391, 'Cultivar IV Región', 'Segunda', 100, 15000, 16000, 15500, '$/bandeja 18 kilos', 'Provincia de Limarí', 861, 18)
    17 = @(44405, 'Cultivar IV Región', 'Segunda', 140, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18)
    18 = @(44398, 'Cultivar IV Región', 'Primera', 100, 17000, 18000, 17500, '$/bandeja 18 kilos', 'Provincia de Limarí', 972, 18)
    19 = @(44398, 'Cultivar IV Región', 'Segunda', 100, 15000, 16000, 15500, '$/bandeja 18 kilos', 'Provincia de Limarí', 861, 18)
    20 = @(44454, 'Cultivar IV Región', 'Primera', 160, 19000, 20000, 19500, '$/bandeja 18 kilos', 'Provincia de Limarí', 1083, 18)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D: Fecha
    $ws.Cells.Item($r, 8).Value  = $vals[1]   # H: Variedad
    $ws.Cells.Item($r, 9).Value  = $vals[2]   # I: Calidad
    $ws.Cells.Item($r, 10).Value = $vals[3]   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $vals[4]   # K: Precio mínimo
    $ws.Cells.Item($r, 12).Value = $vals[5]   # L: Precio máximo
    $ws.Cells.Item($r, 13).Value = $vals[6]   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $vals[7]   # N: Unidad de comercialización
    $ws.Cells.Item($r, 15).Value = $vals[8]   # O: Origen
    $ws.Cells.Item($r, 16).Value = $vals[9]   # P: Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $vals[10]  # Q: Kg o Unidades
}
